$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the URL (matchsource -> matchsync)
$ws.Range("B2").Value = "http://fhir.nmdp.org/ig/matchsync/ValueSet/nmdp-race-codes"

# Set the Experimental value to the literal text "true" (was blank).
# A plain Value assignment of "true" is auto-coerced to a Boolean by Excel,
# so instead compute it as a formula and paste back as a value so the
# cell ends up holding literal text (keeps the existing cell style too).
$cell = $ws.Range("B7")
$cell.Formula = "=""true"""
$cell.Copy()
$cell.PasteSpecial(-4163)

# Update the Date value
$ws.Range("B8").Value = "2024-02-19T18:37:26-06:00"
